# Adds a new data-collection wave (11. 11. 2021) to both worksheets:
#  - "data"   sheet: new column AJ (col 36), header + percentage values rows 2-61,
#              and updates the title row (62) date text.
#  - "pocetR" sheet: new column AI (col 35), header + sample-size values rows 2-24,
#              and updates the title row (25) date text plus an empty tail cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "data"  (percentages)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# New header cell AJ1, copying the formatting (bold/border/centered) of AI1.
$ws1.Cells.Item(1, 36).Value = "11. 11. 2021"
$ws1.Cells.Item(1, 35).Copy()
$ws1.Cells.Item(1, 36).PasteSpecial(-4122)

$sheet1Values = @(
  0.53,0.31,0.16,0.34,0.21,0.45,0.55,0.32,0.13,0.55,
  0.29,0.16,0.48,0.31,0.21,0.5600000000000001,0.29,0.15,0.46,0.31,
  0.23,0.49,0.37,0.14,0.43,0.35,0.22,0.57,0.28,0.15,
  0.61,0.28,0.11,0.42,0.32,0.26,0.39,0.38,0.23,0.64,
  0.24,0.12,0.68,0.23,0.09,0.51,0.29,0.2,0.7,0.26,
  0.04,0.45,0.33,0.22,0.68,0.17,0.15,0.61,0.29,0.1
)

$r = 2
foreach ($v in $sheet1Values) {
  $ws1.Cells.Item($r, 36).Value = $v
  $r = $r + 1
}

# Title / footer row (row 62), update the "aktualizace" date in the text.
$ws1.Cells.Item(62, 1).Value = "Život během pandemie, Obavy ze ztráty práce, % respondentů celkově a ve skupinách, aktualizace 18. 11. 2021"

# ---------------------------------------------------------------------------
# Sheet 2: "pocetR"  (sample sizes)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# New header cell AI1, copying the formatting of AH1.
$ws2.Cells.Item(1, 35).Value = "11. 11. 2021"
$ws2.Cells.Item(1, 34).Copy()
$ws2.Cells.Item(1, 35).PasteSpecial(-4122)

$sheet2Values = @(
  955,88,867,744,145,6,60,711,133,59,
  52,353,373,229,108,295,301,151,258,90,
  254,134,83
)

$r = 2
foreach ($v in $sheet2Values) {
  $ws2.Cells.Item($r, 35).Value = $v
  $r = $r + 1
}

# Title row (row 25): update date text and add the trailing empty cell (AI25)
# matching the blank inlineStr cells already present across the row. A plain
# empty-string assignment is treated as a no-op by the engine (the cell is
# never materialised), so force the cell to exist by touching its formatting
# and then clearing that formatting again, which leaves a present-but-empty
# cell behind.
$ws2.Cells.Item(25, 1).Value = "Život během pandemie, Obavy ze ztráty práce, velikost dotázaného souboru celkově a ve skupinách, aktualizace 18. 11. 2021"
$ws2.Cells.Item(25, 35).NumberFormat = "General"
$ws2.Cells.Item(25, 35).ClearFormats()
